# SonoVerse.xlsx update
# - Remove the "Biliary hamartomatosis" row (old row 4)
# - Add a new "Biliary Hamartomas (von Meyenburg Complexes)" row for Liver
#   (ends up sorted at row 13), with a plain-text YouTube link (not a live
#   hyperlink) in column D, same as a few other D-cells in this sheet.
# - All the other rows shift up by one to close the gap left by the deleted
#   row, keeping the A-column alphabetical sort state intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for rows 1-24, columns A-D (column E has no data in this sheet).
$rows = @{
    1  = @("Organ", "Term", "Variant/Detail", "YouTube Link")
    2  = @("Bowel", "Coming soon", "", "")
    3  = @("Gallbladder and biliary tract", "Adenomyomatosis ", "Clip 1 B-mode", "https://youtu.be/zxTC0YBY2RY ")
    4  = @("Gallbladder and biliary tract", "Intrahepatic bile ducts mild dilatation", "Clip 1 B-mode", "")
    5  = @("Kidney", "Renal stone", "Clip 1 B-mode + Color", "https://youtu.be/2kRZcpi70Aw ")
    6  = @("Urinary tract", "Ureteral jets", "Clip 1 B-mode", "")
    7  = @("Liver", "Metastasis target-like", "Clip 1 B-mode ", "")
    8  = @("Liver", "HNF1α-mutated hepatocellular adenoma - Hyperechoic", "Clip 1 B-mode + Color ", "https://youtu.be/91M82AIMyu0 ")
    9  = @("Liver", "Focal nodular hyperplasia - Isoechoic", "Clip 1 B-mode + Color + microV", "https://youtu.be/rg0sFcu0rVQ")
    10 = @("Liver", "Hepatocellular carcinoma - Heterogeneous nodular lesions ", "Clip 1 B-mode + Color + microV", "https://youtu.be/15o_Km86IzM ")
    11 = @("Liver", "Perihepatic reactive lymph nodes", "Clip 1-Bmode", "https://youtu.be/kaROVVBl9Bc")
    12 = @("Liver", "Cavernous Hemangioma - Iso-Hyperechoic, Trilobulated", "Clip 1 B-mode + Color", "https://youtu.be/RhSUFLTmTl4")
    13 = @("Liver", "Biliary Hamartomas (von Meyenburg Complexes) ", "Clip 1 B-mode", "https://youtu.be/knQhiK4Y7kY")
    14 = @("Liver vasculature", "Portal vein thrombosis", "Clip 1 B-mode + Color", "https://youtu.be/DjI1kEnzfSQ ")
    15 = @("Liver vasculature", "Spontaneous intrahepatic porto-systemic shunt", "Clip 1 B-mode + Color", "https://youtu.be/U3ydTsRwxok ")
    16 = @("Liver vasculature", "Congestive Hepatopathy ", "Clip 1-Bmode + Color Doppler", "https://youtu.be/sRu_NTopG3Y")
    17 = @("Lymph nodes", "Coming soon", "", "")
    18 = @("Miscellaneous", "Adrenal adenoma", "Clip 1 B-mode", "https://youtu.be/xBfd04F4Ni8 ")
    19 = @("Pancreas", "Acute necrotizing pancreatitis ", "Clip 1 B-mode", "https://youtu.be/JvwODCASLYQ ")
    20 = @("Pancreas", "Neuroendocrine Tumor G1 – Hypoechoic", "Clip 1 B-mode + Color + microV", "https://youtu.be/pc-vbxSRTbs ")
    21 = @("Pancreas", "Stones in the Main Pancreatic Duct (Pancreatolithiasis)", "Clip 1 B-mode + Color", "https://youtu.be/Axbee4vjNtU")
    22 = @("Spleen", "Splenic calcification with posterior shadowing", "Clip 1 B-mode", "https://youtu.be/qushjTAy6XQ ")
    23 = @("Spleen", "Accessory spleen", "Clip 1 B-mode", "https://youtu.be/_FckFwJwynI ")
    24 = @("Thyroid", "Isoechoic nodule with peripheral calcifications", "Clip 1 B-mode + Color", "https://youtu.be/z_oaRVxRz5s ")
}

# Drop every existing hyperlink relationship first; the correct ones (and
# only those) get re-added below, after all the text is back in place.
$ws.Hyperlinks.Delete()

# Rewrite every data row's A:D cells so the sheet matches the new layout.
for ($r = 1; $r -le 24; $r++) {
    $vals = $rows[$r]
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Column D used the "hyperlink look" style on every filled cell (including a
# few cells that are plain text, not live links) - reapply that to every row
# with a YouTube-link value (including the newly inserted row 13), and make
# sure every empty D cell is back to the plain default style.
$linkRows = @(3, 5, 8, 9, 10, 11, 12, 13, 14, 15, 16, 18, 19, 20, 21, 22, 23, 24)
$plainDRows = @(1, 2, 4, 6, 7, 17)
foreach ($r in $linkRows) {
    $ws.Range("D" + $r).Style = "Collegamento ipertestuale"
}
foreach ($r in $plainDRows) {
    $ws.Range("D" + $r).Style = "Normale"
}

# B5 ("Renal stone") carries the vertical-center alignment style that used
# to sit on the pre-edit "Renal stone" row.
$ws.Range("B5").Style = "Normale"
$ws.Range("B5").VerticalAlignment = -4108

# Re-create live hyperlinks at their new positions. D13 (the new Biliary
# Hamartomas row) intentionally stays plain text with no live link, matching
# the source workbook.
$ws.Hyperlinks.Add($ws.Range("D3"), "https://youtu.be/zxTC0YBY2RY")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://youtu.be/2kRZcpi70Aw")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://youtu.be/91M82AIMyu0")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/RhSUFLTmTl4")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://youtu.be/DjI1kEnzfSQ")
$ws.Hyperlinks.Add($ws.Range("D15"), "https://youtu.be/U3ydTsRwxok")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/15o_Km86IzM")
$ws.Hyperlinks.Add($ws.Range("D18"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D19"), "https://youtu.be/JvwODCASLYQ")
$ws.Hyperlinks.Add($ws.Range("D20"), "https://youtu.be/pc-vbxSRTbs")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://youtu.be/Axbee4vjNtU")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://youtu.be/qushjTAy6XQ")
$ws.Hyperlinks.Add($ws.Range("D23"), "https://youtu.be/_FckFwJwynI")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/z_oaRVxRz5s")

# Match the saved selection from the edit.
$ws.Range("D16").Select()
